$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SKILL")

# Insert a new column at H, shifting the old "Test Field"/Annotation column to I.
$ws.Columns("H").Insert()

# Undo-then-redo the header cell move so the vacated H1 doesn't linger as a
# phantom empty cell in the exported XML: shift the old value back out of I1
# into H1, clear H1 completely (content + format), then give I1 its real,
# final value as a clean write.
$ws.Range("H1").Delete(-4159)
$ws.Range("I1").Value = "Annotation"
$ws.Range("H1").Clear()

# Populate the new "cost" / CBigInt field column.
$ws.Range("H3").Value = "CBigInt"
$ws.Range("H2").Value = "cost"
$ws.Range("H5").Value = "1.2, 30"
$ws.Range("H6").Value = "1.1, 50"
$ws.Range("H4").Value = "class"

$ws.Columns("H").ColumnWidth = 12.75
$ws.Columns("I").AutoFit() | Out-Null

$ws.Range("H1").Select()
